# Add "Ketamin" and "Midazolam" dose-calculator tabs to the workbook.
# Both new sheets reuse the exact row/column/style layout of the existing
# "Lidocaine" sheet (same formulas, same style ids), so we duplicate that
# sheet twice and then just swap in the new labels/values.
#
# Note: worksheet object handles returned by Item(index) are positional,
# so once all the inserts/renames are done we re-fetch every sheet we still
# need **by name** before touching it again.

$wb = $excel.ActiveWorkbook
$lidocaine = $wb.Worksheets.Item("Lidocaine")

# --- Create "Midazolam" first (copied right after Lidocaine) -------------
$lidocaine.Copy($null, $lidocaine)
$wb.Worksheets.Item($lidocaine.Index + 1).Name = "Midazolam"

# --- Create "Ketamin" (copied right after Lidocaine, i.e. before Midazolam)
$lidocaine.Copy($null, $lidocaine)
$wb.Worksheets.Item($lidocaine.Index + 1).Name = "Ketamin"

# Re-fetch fresh, name-bound references now that the tab order is final.
$ketamin = $wb.Worksheets.Item("Ketamin")
$midazolam = $wb.Worksheets.Item("Midazolam")

# ---------------------------------------------------------------------
# Ketamin tab values
# ---------------------------------------------------------------------
$ketamin.Range("A1").Value = "Ketamin 150mg/kg "
$ketamin.Range("D3").Value = 150
$ketamin.Range("D5").Value = 36
$ketamin.Range("C8").Value = "Concentration of Ketamin"
$ketamin.Range("D8").Value = 100
$ketamin.Range("D9").Value = 0.2
$ketamin.Range("C13").Value = ""

# ---------------------------------------------------------------------
# Midazolam tab values
# ---------------------------------------------------------------------
$midazolam.Range("A1").Value = "Midazolam  1.5mg/kg "
$midazolam.Range("D3").Value = 1.5
$midazolam.Range("D5").Value = 36
$midazolam.Range("C8").Value = "Concentration of Midazolam"
$midazolam.Range("D8").Value = 5
$midazolam.Range("D9").Value = 0.04
$midazolam.Range("C13").Value = ""

# ---------------------------------------------------------------------
# View niceties matching the authored workbook (zoom + final selections)
# ---------------------------------------------------------------------
$meloxicam = $wb.Worksheets.Item("Meloxicam")
$pentobarbital = $wb.Worksheets.Item("Pentobarbital")

$meloxicam.Activate()
$excel.ActiveWindow.Zoom = 130

$pentobarbital.Activate()
$excel.ActiveWindow.Zoom = 130

$lidocaine = $wb.Worksheets.Item("Lidocaine")
$lidocaine.Activate()
$excel.ActiveWindow.Zoom = 130
$lidocaine.Range("D10").Select()

$ketamin = $wb.Worksheets.Item("Ketamin")
$ketamin.Activate()
$excel.ActiveWindow.Zoom = 130
$ketamin.Range("C13").Select()

$midazolam = $wb.Worksheets.Item("Midazolam")
$midazolam.Activate()
$excel.ActiveWindow.Zoom = 130
$midazolam.Range("D11").Select()
